$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(22, -61.64058685302734, 21.88947105407715, 0.6875860691070557),
    @(23, -27.69689178466797, -7.563150405883789, -14.80709171295166),
    @(24, -1.68331241607666, -72.37289428710938, 70.38520812988281),
    @(25, -11.52985954284668, 10.27250671386719, 4.496110916137695),
    @(26, 4.430462837219238, -20.51617813110352, 37.68183135986328),
    @(27, -80.67318725585938, -3.729990005493164, -17.92157745361328),
    @(28, -7.960968494415283, 2.627067565917969, -1.957346677780152),
    @(29, -12.41016101837158, -18.66624450683594, 3.092851161956787),
    @(30, 6.215286254882812, 5.965664863586426, 14.03029632568359),
    @(31, -6.658430576324463, 14.9360237121582, 14.47612380981445)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
}
